# "add Results_Alex + edit final results"
#
# The "Gesamt" (C) column for each recommender is the average of the three
# raters' scores (Peter = D, Alex = E, Max = F). Alex's ("Results_Alex")
# scores were missing for most rows; this adds them in and lets the
# "Gesamt" averages recompute. Row 11 ("Meta-Mix") never got an Alex score,
# so its formula is rewritten to just average Peter (D11) and Max (F11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New Alex (column E) ratings
$ws.Range("E9").Value  = 5.4
$ws.Range("E10").Value = 6.15
$ws.Range("E12").Value = 6.08
$ws.Range("E13").Value = 6.2
$ws.Range("E14").Value = 5.93

# "Meta-Mix" row has no Alex rating -> average only Peter & Max directly
# instead of the D10:C14 shared "(D+E+F)/3" formula
$ws.Range("C11").Formula = "=(F11+D11)/2"

# The other rows keep their existing "(Dn+En+Fn)/3" shared formula; it will
# simply pick up the newly-entered Alex values on recalculation.
$ws.Calculate()

# Move the selection/viewport: previously scrolled to A7 with G20 selected,
# now back at the top with G10 selected
$ws.Range("G10").Select()
